$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(9)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Step 1: rebuild the text frame's paragraphs from scratch (same text, plus one
# extra trailing blank paragraph) so PowerPoint regenerates clean run
# properties for every paragraph - this clears the stale err="1" spell-flags
# and the superfluous endParaRPr left over on paragraphs 2 and 3, without
# altering the visible content yet.
$tr.Text = "Manque de temps`rDddd`rffff`r`r"
$tr2 = $tf.TextRange
$tr2.Paragraphs(5).Delete()

# Step 2: apply the real text + font-size edits per paragraph so the existing
# (now clean) paragraph 4 is left completely untouched.
$tr3 = $tf.TextRange
$tr3.Paragraphs(1).Text = "Manque de temps"
$tr3.Paragraphs(1).Font.Size = 32

$tr3.Paragraphs(2).Text = "Un commercial trop ambitieux"
$tr3.Paragraphs(2).Font.Size = 32

$tr3.Paragraphs(3).Text = "Mauvaise estimation de la complexité de sortir un produit «jouable»"
$tr3.Paragraphs(3).Font.Size = 32

$tr3.Paragraphs(4).Font.Size = 20
